$d = $word.ActiveDocument

# --- Step 1: the old "Development steps:" heading becomes "Traffic Control Env" ---
$devParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "Development steps:") {
        $devParaIndex = $i
        break
    }
}
$anchorPara = $d.Paragraphs.Item($devParaIndex)
$anchorPara.Range.Find.Execute("Development steps:", $true, $false, $false, $false, $false, $true, 1, $false, "Traffic Control Env", 2)

# --- Step 2: insert 5 new paragraphs right after it, anchored off that same range ---
$anchorRange = $anchorPara.Range
for ($i = 0; $i -lt 5; $i++) {
    $anchorRange.InsertParagraphAfter()
}

# The 5 freshly-created (still empty) paragraphs now sit immediately after $devParaIndex
$sumoPara        = $d.Paragraphs.Item($devParaIndex + 1)
$matplotlibPara  = $d.Paragraphs.Item($devParaIndex + 2)
$blankPara1      = $d.Paragraphs.Item($devParaIndex + 3)
$blankPara2      = $d.Paragraphs.Item($devParaIndex + 4)
$newDevStepsPara = $d.Paragraphs.Item($devParaIndex + 5)

$sumoPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Needs the SUMO package (under scripts, ubuntu_setup.sh)</w:t></w:r></w:p>')

$matplotlibPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Plus some other regular things like matplotlib</w:t></w:r></w:p>')

$blankPara1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>')

$blankPara2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>')

$newDevStepsPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Development steps:</w:t></w:r></w:p>')

Write-Output "done"
